# Apply the committed changes to the "Form_Responses" table:
#  1. Simplify the rich-text header in L1 ("ส่วน <bold+underline>" + "<bold>...")
#     down to plain bold text "ความพึงพอใจต่อบริการของโรงพยาบาลในภาพรวม"
#     (this also renames the corresponding table column automatically).
#  2. Append three new survey responses as rows 225-227 of the
#     Form_Responses table, copying the row-224 cell formatting so the new
#     rows look like the existing ones, then writing in the actual values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# --- 1. Fix the section header in L1 -------------------------------------
$ws.Range("L1").Value = "ความพึงพอใจต่อบริการของโรงพยาบาลในภาพรวม"
$ws.Range("L1").Font.Bold = $true
$ws.Range("L1").Font.Underline = $false

# --- 2. Append the three new rows -----------------------------------------

# Row 225
$lo.ListRows.Add() | Out-Null
$ws.Range("A224:D224").Copy($ws.Range("A225:D225"))
$ws.Range("G224:Y224").Copy($ws.Range("G225:Y225"))

$ws.Range("A225").Value = 45957.43552460648
$ws.Range("B225").Value = "หน่วยตรวจหู คอ จมูก"
$ws.Range("C225").Value = "ครั้งแรก"
$ws.Range("D225").Value = "สุขภาพไม่ดี"
$ws.Range("G225").Value = "อายุ 52 - 70 ปี"
$ws.Range("H225").Value = "ภายในจังหวัดเชียงราย"
$ws.Range("I225").Value = "รับจ้าง"
$ws.Range("J225").Value = "จ่ายตรง"
$ws.Range("K225").Value = 244284.0
$ws.Range("L225").Value = "มากที่สุด"
$ws.Range("M225").Value = "สะดวกมาก"
$ws.Range("N225").Value = "สะดวกมาก"
$ws.Range("O225").Value = "เหมาะสมมาก"
$ws.Range("P225").Value = "มากที่สุด"
$ws.Range("Q225").Value = "ชัดเจนมาก"
$ws.Range("R225").Value = "มากที่สุด"
$ws.Range("S225").Value = "เหมาะสมมาก"
$ws.Range("T225").Value = "ชัดเจนมาก"
$ws.Range("U225").Value = "ชัดเจนมาก"
$ws.Range("V225").Value = "ชัดเจนมาก"
$ws.Range("W225").Value = "ใช่"
$ws.Range("X225").Value = "ใช่"
$ws.Range("Y225").Value = "ไม่มี"

# Row 226
$lo.ListRows.Add() | Out-Null
$ws.Range("A224:E224").Copy($ws.Range("A226:E226"))
$ws.Range("G224:Y224").Copy($ws.Range("G226:Y226"))
$ws.Range("AA224").Copy($ws.Range("AA226"))

$ws.Range("A226").Value = 45957.462592881944
$ws.Range("B226").Value = "หน่วยแพทย์บูรณาการ"
$ws.Range("C226").Value = "มากกว่า 1 ครั้ง"
$ws.Range("D226").Value = "สุขภาพดี"
$ws.Range("E226").Value = "มีอาการปวดเมื่อย"
$ws.Range("G226").Value = "อายุ 36 - 51 ปี"
$ws.Range("H226").Value = "พะเยา"
$ws.Range("I226").Value = "รับราชการ"
$ws.Range("J226").Value = "กรมบัญชีกลาง"
$ws.Range("K226").Value = 45957.0
$ws.Range("L226").Value = "มากที่สุด"
$ws.Range("M226").Value = "สะดวกมาก"
$ws.Range("N226").Value = "สะดวกมาก"
$ws.Range("O226").Value = "เหมาะสมมาก"
$ws.Range("P226").Value = "มากที่สุด"
$ws.Range("Q226").Value = "ชัดเจนมาก"
$ws.Range("R226").Value = "มากที่สุด"
$ws.Range("S226").Value = "เหมาะสมมาก"
$ws.Range("T226").Value = "ชัดเจนมาก"
$ws.Range("U226").Value = "ชัดเจนมาก"
$ws.Range("V226").Value = "ชัดเจนมาก"
$ws.Range("W226").Value = "ใช่"
$ws.Range("X226").Value = "ใช่"
$ws.Range("Y226").Value = "ไม่มี"
$ws.Range("AA226").Value = "ทำดีแล้วครับ รักษามาตรฐานไว้ต่อไป พัฒนาขึ้นเรื่อยๆ"

# Row 227
$lo.ListRows.Add() | Out-Null
$ws.Range("A224:E224").Copy($ws.Range("A227:E227"))
$ws.Range("G224:Y224").Copy($ws.Range("G227:Y227"))

$ws.Range("A227").Value = 45957.47176168981
$ws.Range("B227").Value = "หน่วยตรวจเด็กสุขภาพดี"
$ws.Range("C227").Value = "มากกว่า 1 ครั้ง"
$ws.Range("D227").Value = "สุขภาพดี"
$ws.Range("E227").Value = "วัคซีน"
$ws.Range("G227").Value = "อายุ 18 - 35 ปี"
$ws.Range("H227").Value = "ภายในจังหวัดเชียงราย"
$ws.Range("I227").Value = "กฟส.แม่สาย"
$ws.Range("J227").Value = "รัฐวิสาหกิจ"
$ws.Range("K227").Value = 45957.0
$ws.Range("L227").Value = "มากที่สุด"
$ws.Range("M227").Value = "สะดวกมาก"
$ws.Range("N227").Value = "สะดวกมาก"
$ws.Range("O227").Value = "เหมาะสมมาก"
$ws.Range("P227").Value = "มากที่สุด"
$ws.Range("Q227").Value = "ชัดเจนมาก"
$ws.Range("R227").Value = "มากที่สุด"
$ws.Range("S227").Value = "เหมาะสมมาก"
$ws.Range("T227").Value = "ชัดเจนมาก"
$ws.Range("U227").Value = "ชัดเจนมาก"
$ws.Range("V227").Value = "ชัดเจนมาก"
$ws.Range("W227").Value = "ใช่"
$ws.Range("X227").Value = "ใช่"
$ws.Range("Y227").Value = "ไม่มี"

Write-Host "Edit applied"
